$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 126.14286
$ws.Range("I9").Value = 126.14286
$ws.Range("K9").Value = 126.14286
$ws.Range("M9").Value = 42.85714

$ws.Range("H32").Value = 1072.2222
$ws.Range("I32").Value = 764.5
$ws.Range("J32").Value = 1318.4
$ws.Range("K32").Value = 764.5
$ws.Range("L32").Value = 1318.4
$ws.Range("M32").Value = -438.5
$ws.Range("N32").Value = -1970.4

$ws.Range("H62").Value = 2424.0908
$ws.Range("I62").Value = 2197.8572
$ws.Range("J62").Value = 2820
$ws.Range("K62").Value = 2197.8572
$ws.Range("L62").Value = 2820
$ws.Range("M62").Value = -1573.8572
$ws.Range("N62").Value = -4068

$ws.Range("H64").Value = 2870.0625
$ws.Range("I64").Value = 2302.8
$ws.Range("J64").Value = 3370.5881
$ws.Range("K64").Value = 2302.8
$ws.Range("L64").Value = 3370.5881
$ws.Range("M64").Value = -2054.8
$ws.Range("N64").Value = -3866.5881

$ws.Range("H65").Value = 2424.0908
$ws.Range("I65").Value = 2197.8572
$ws.Range("J65").Value = 2820
$ws.Range("K65").Value = 10989.286
$ws.Range("L65").Value = 14100
$ws.Range("M65").Value = -7869.286
$ws.Range("N65").Value = -20340

$ws.Range("H67").Value = 2870.0625
$ws.Range("I67").Value = 2302.8
$ws.Range("J67").Value = 3370.5881
$ws.Range("K67").Value = 2302.8
$ws.Range("L67").Value = 3370.5881
$ws.Range("M67").Value = -1444.8
$ws.Range("N67").Value = -5086.5881

$ws.Range("H80").Value = 6467
$ws.Range("I80").Value = 641.4286
$ws.Range("J80").Value = 9865.25
$ws.Range("K80").Value = 1924.2858
$ws.Range("L80").Value = 29595.75
$ws.Range("M80").Value = -926.2857999999999
$ws.Range("N80").Value = -31591.75

$ws.Range("H83").Value = 6467
$ws.Range("I83").Value = 641.4286
$ws.Range("J83").Value = 9865.25
$ws.Range("K83").Value = 5772.8574
$ws.Range("L83").Value = 88787.25
$ws.Range("M83").Value = -780.8573999999999
$ws.Range("N83").Value = -98771.25

$ws.Range("H100").Value = 3151.375
$ws.Range("I100").Value = 1200
$ws.Range("J100").Value = 3430.1428
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 3430.1428
$ws.Range("M100").Value = -659
$ws.Range("N100").Value = -4512.1428

$ws.Range("H129").Value = 954.38
$ws.Range("I129").Value = 499.69232
$ws.Range("J129").Value = 1114.1351
$ws.Range("K129").Value = 1499.07696
$ws.Range("L129").Value = 3342.4053
$ws.Range("M129").Value = 3500.92304
$ws.Range("N129").Value = -13342.4053

$ws.Range("H132").Value = 3533.2222
$ws.Range("I132").Value = 3340.2104
$ws.Range("J132").Value = 3748.9412
$ws.Range("K132").Value = 10020.6312
$ws.Range("L132").Value = 11246.8236
$ws.Range("M132").Value = -7490.6312
$ws.Range("N132").Value = -16306.8236

$ws.Range("H137").Value = 1532.8529
$ws.Range("I137").Value = 1236.7727
$ws.Range("J137").Value = 2075.6667
$ws.Range("K137").Value = 3710.3181
$ws.Range("L137").Value = 6227.000100000001
$ws.Range("M137").Value = -1160.3181
$ws.Range("N137").Value = -11327.0001

$ws.Range("H138").Value = 3494.9856
$ws.Range("I138").Value = 2577.75
$ws.Range("J138").Value = 3688.0876
$ws.Range("K138").Value = 7733.25
$ws.Range("L138").Value = 11064.2628
$ws.Range("M138").Value = -2593.25
$ws.Range("N138").Value = -21344.2628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 625.5714
$ws.Range("I4").Value = 635.2308
$ws.Range("K4").Value = 635.2308
$ws.Range("M4").Value = -519.2308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 20007
$ws.Range("J15").Value = 20007
$ws.Range("L15").Value = 20007
$ws.Range("N15").Value = -20461

$ws.Range("H22").Value = 17330
$ws.Range("I22").Value = 17330
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 17330
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -17157
$ws.Range("N22").ClearContents()

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4664

$ws.Range("H36").Value = 1189.75
$ws.Range("I36").Value = 1189.75
$ws.Range("K36").Value = 1189.75
$ws.Range("M36").Value = -655.75

$ws.Range("H105").Value = 4119.3335
$ws.Range("I105").Value = 3877
$ws.Range("K105").Value = 3877
$ws.Range("M105").Value = -2130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 766.6667
$ws.Range("I105").Value = 766.6667
$ws.Range("K105").Value = 766.6667
$ws.Range("M105").Value = 980.3333

$ws.Range("H107").Value = 409.58334
$ws.Range("I107").Value = 507.125
$ws.Range("K107").Value = 507.125
$ws.Range("M107").Value = 1412.875

$ws.Range("H123").Value = 63780
$ws.Range("J123").Value = 63780
$ws.Range("L123").Value = 63780
$ws.Range("N123").Value = -73580

$ws.Range("H132").Value = 1743.3182
$ws.Range("J132").Value = 4332.6665
$ws.Range("L132").Value = 12997.9995
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1595.5714
$ws.Range("I10").Value = 296.66666
$ws.Range("J10").Value = 2569.75
$ws.Range("K10").Value = 889.9999799999999
$ws.Range("L10").Value = 7709.25
$ws.Range("M10").Value = -750.9999799999999
$ws.Range("N10").Value = -7987.25

$ws.Range("H131").Value = 38463730
$ws.Range("J131").Value = 41668990
$ws.Range("L131").Value = 125006970
$ws.Range("N131").Value = -125017050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 30000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 30000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -29732
$ws.Range("N34").ClearContents()

$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 30000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -29685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 30000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -28908
$ws.Range("N79").ClearContents()

$ws.Range("H102").Value = 4059.7778
$ws.Range("I102").Value = 3804.8
$ws.Range("J102").Value = 4378.5
$ws.Range("K102").Value = 3804.8
$ws.Range("L102").Value = 4378.5
$ws.Range("M102").Value = -2182.8
$ws.Range("N102").Value = -7622.5

$ws.Range("H122").Value = 1987.4642
$ws.Range("I122").Value = 2283.1052
$ws.Range("J122").Value = 1363.3334
$ws.Range("K122").Value = 6849.3156
$ws.Range("L122").Value = 4090.0002
$ws.Range("M122").Value = -4399.3156
$ws.Range("N122").Value = -8990.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6141.0356
$ws.Range("I2").Value = 955.6667
$ws.Range("J2").Value = 8597.263000000001
$ws.Range("K2").Value = 955.6667
$ws.Range("L2").Value = 8597.263000000001
$ws.Range("M2").Value = -843.6667
$ws.Range("N2").Value = -8821.263000000001

$ws.Range("H12").Value = 15001.5
$ws.Range("I12").Value = 10003
$ws.Range("J12").Value = 20000
$ws.Range("K12").Value = 10003
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = -9833
$ws.Range("N12").Value = -20340

$ws.Range("H46").Value = 1733.2609
$ws.Range("I46").Value = 1682.5
$ws.Range("J46").Value = 1849.2858
$ws.Range("K46").Value = 1682.5
$ws.Range("L46").Value = 1849.2858
$ws.Range("M46").Value = -1494.5
$ws.Range("N46").Value = -2225.2858

$ws.Range("H94").Value = 18063.334
$ws.Range("J94").Value = 18063.334
$ws.Range("L94").Value = 18063.334
$ws.Range("N94").Value = -19415.334

$ws.Range("H132").Value = 4180.2095
$ws.Range("I132").Value = 3420.4827
$ws.Range("J132").Value = 5753.9287
$ws.Range("K132").Value = 10261.4481
$ws.Range("L132").Value = 17261.7861
$ws.Range("M132").Value = -7731.4481
$ws.Range("N132").Value = -22321.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6000
$ws.Range("J2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("N2").Value = -6224

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
